$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Kelp consumption")

# --- Data updates on "Kelp consumption" sheet ---

# New starvation-time related values added to row 2 (columns N, O)
$ws2.Range("N2").Value = 1900
$ws2.Range("O2").Value = 700

# Correct the trial date (col B) and Julian date (col C) for rows 58-71
# (these were recorded with the wrong date and are shifted back by 3 days)
for ($r = 58; $r -le 71; $r++) {
    $ws2.Cells.Item($r, 2).Value = 43713
    $ws2.Cells.Item($r, 3).Value = 248
}

# Correct the trial date (col B) and Julian date (col C) for rows 72-77
for ($r = 72; $r -le 77; $r++) {
    $ws2.Cells.Item($r, 2).Value = 43714
    $ws2.Cells.Item($r, 3).Value = 249
}

# A new (blank) row was added below the data, with column B formatted
# the same way (date format) as the rest of the column
$ws2.Cells.Item(78, 2).NumberFormat = $ws2.Cells.Item(77, 2).NumberFormat

# --- View state updates ---

# "Kelp consumption" tab becomes the active/selected tab, scrolled back
# to the top of the frozen pane, with H63 as the active selection
$ws2.Activate() | Out-Null
$ws2.Range("A2").Select() | Out-Null
$ws2.Range("H63").Select() | Out-Null
